$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# ---- New skill rows (201 火球术, 202 冰弹术, 203 闪电术) ----
$newSkills = @(
    @{ Row = 23; Id = 201; Name = "火球术"; Desc = "在那个时代人人都会的火球术" },
    @{ Row = 24; Id = 202; Name = "冰弹术"; Desc = "在那个时代人人都会的冰弹术" },
    @{ Row = 25; Id = 203; Name = "闪电术"; Desc = "在那个时代人人都会的闪电术" }
)

# --- Formatting pass -------------------------------------------------
# Register the two "plain" styles (default font for the description
# column, name font for the skill column) before any text is typed, so
# the new cellXfs entries land in the same order as in the target file.
foreach ($skill in $newSkills) {
    $r = $skill.Row
    $ws.Cells.Item($r, 10).VerticalAlignment = -4107
}

foreach ($skill in $newSkills) {
    $r = $skill.Row
    $ws.Range("B3").Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).VerticalAlignment = -4107

    $ws.Range("C3").Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false

# --- Data pass ---------------------------------------------------------
# Skill id + numeric columns.
foreach ($skill in $newSkills) {
    $r = $skill.Row
    $ws.Cells.Item($r, 1).Value = $skill.Id
    $ws.Cells.Item($r, 3).Value = 1
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 1
    $ws.Cells.Item($r, 7).Value = 1
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 9).Value = 0
}

# Skill names first (column B), then descriptions (column J) - matches
# the interning order of the new shared-string entries.
foreach ($skill in $newSkills) {
    $ws.Cells.Item($skill.Row, 2).Value = $skill.Name
}
foreach ($skill in $newSkills) {
    $ws.Cells.Item($skill.Row, 10).Value = $skill.Desc
}

# Matches the author's final cursor position after entering the new rows.
$ws.Range("J25").Select() | Out-Null
